$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: Create Transaction Controller
$ws.Range("A15").Value = 43504
$ws.Range("A15").NumberFormat = "DD/MM/YY"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Implementation"
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = "Create Transaction Controller"
$ws.Range("F15").Value = "Requirements – don’t forget they’re up there and matter!"

# Row 16: Create Merchant Controller
$ws.Range("A16").Value = 43504
$ws.Range("A16").NumberFormat = "DD/MM/YY"
$ws.Range("B16").Value = 0.25
$ws.Range("C16").Value = "Implementation"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = "Create Merchant Controller"
$ws.Range("F16").Value = "basically a copy of above, tested."

$ws.Range("B7").Select() | Out-Null
